$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 164-165, pushing existing rows 164-267 down to 166-269.
$ws.Range("A164:R165").EntireRow.Insert()

# New row 164 ("Primera" quality record dated 2021-10-19)
$ws.Cells.Item(164, 1).Value = 9
$ws.Cells.Item(164, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(164, 3).Value = "Metropolitana"
$ws.Cells.Item(164, 4).Value = 44488
$ws.Cells.Item(164, 5).Value = 13
$ws.Cells.Item(164, 6).Value = 100114014
$ws.Cells.Item(164, 7).Value = "Betarraga"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 5200
$ws.Cells.Item(164, 11).Value = 100
$ws.Cells.Item(164, 12).Value = 110
$ws.Cells.Item(164, 13).Value = 105
$ws.Cells.Item(164, 14).Value = "`$/unidad"
$ws.Cells.Item(164, 15).Value = "Región Metropolitana"
$ws.Cells.Item(164, 16).Value = 105
$ws.Cells.Item(164, 17).Value = 1
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# New row 165 ("Segunda" quality record dated 2021-10-19)
$ws.Cells.Item(165, 1).Value = 9
$ws.Cells.Item(165, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(165, 3).Value = "Metropolitana"
$ws.Cells.Item(165, 4).Value = 44488
$ws.Cells.Item(165, 5).Value = 13
$ws.Cells.Item(165, 6).Value = 100114014
$ws.Cells.Item(165, 7).Value = "Betarraga"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Segunda"
$ws.Cells.Item(165, 10).Value = 3400
$ws.Cells.Item(165, 11).Value = 70
$ws.Cells.Item(165, 12).Value = 80
$ws.Cells.Item(165, 13).Value = 75
$ws.Cells.Item(165, 14).Value = "`$/unidad"
$ws.Cells.Item(165, 15).Value = "Región Metropolitana"
$ws.Cells.Item(165, 16).Value = 75
$ws.Cells.Item(165, 17).Value = 1
$ws.Cells.Item(165, 18).Value = "Hortaliza"
